# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Prefix with an apostrophe so Excel always stores the value as literal
    # text (never auto-converts look-alike numbers/dates), then reset the
    # cell style so no stray quote-prefix formatting is left behind -
    # matching the plain (unstyled) text cells already in the sheet.
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).Style = "Normal"
}

Set-TextValue "D2" "27.586.05"
Set-TextValue "E2" "  -1.45%  "
Set-TextValue "D3" "1.753.29"
Set-TextValue "E3" "  -0.72%  "
Set-TextValue "E4" "  +0.17%  "
Set-TextValue "D5" "324.42"
Set-TextValue "E5" "  +0.74%  "
Set-TextValue "E6" "  +0.18%  "
Set-TextValue "D7" "0.4588"
Set-TextValue "E7" "  +8.08%  "
Set-TextValue "D8" "0.3591"
Set-TextValue "E8" "  -0.18%  "
Set-TextValue "E9" "  +0.52%  "
Set-TextValue "D10" "42.01"
Set-TextValue "E10" "  -5.11%  "
Set-TextValue "D11" "1.095"
Set-TextValue "E11" "  -0.97%  "
Set-TextValue "D12" "1.001"
Set-TextValue "E12" "  +0.20%  "
Set-TextValue "D13" "20.81"
Set-TextValue "E13" "  -3.04%  "
Set-TextValue "D14" "6.007"
Set-TextValue "E14" "  -1.55%  "
Set-TextValue "D15" "7.099"
Set-TextValue "E15" "  -3.44%  "
Set-TextValue "D16" "1.754.87"
Set-TextValue "E16" "  -2.27%  "
Set-TextValue "D17" "93.32"
Set-TextValue "E17" "  +1.77%  "
Set-TextValue "D18" "0.00001066"
Set-TextValue "E18" "  +0.46%  "
Set-TextValue "D19" "0.06413"
Set-TextValue "E19" "  +0.38%  "
Set-TextValue "D20" "0.9998"
Set-TextValue "E20" "  +0.23%  "
Set-TextValue "D21" "16.80"
Set-TextValue "E21" "  -1.98%  "
Set-TextValue "D22" "5.833"
Set-TextValue "E22" "  -2.50%  "
Set-TextValue "D23" "27.642.08"
Set-TextValue "E23" "  -1.29%  "
Set-TextValue "D24" "11.20"
Set-TextValue "E24" "  -0.91%  "
Set-TextValue "D25" "2.120"
Set-TextValue "E25" "  -1.42%  "
Set-TextValue "D26" "162.66"
Set-TextValue "E26" "  +2.51%  "
Set-TextValue "D27" "20.44"
Set-TextValue "E27" "  +1.23%  "
Set-TextValue "D28" "1.954.34"
Set-TextValue "E28" "  -1.99%  "
Set-TextValue "D29" "2.082"
Set-TextValue "D30" "127.28"
Set-TextValue "E30" "  +0.92%  "
Set-TextValue "D31" "1.082"
Set-TextValue "E31" "  -8.20%  "
Set-TextValue "D32" "0.09156"
Set-TextValue "E32" "  +1.29%  "
Set-TextValue "E33" "  +4.49%  "
Set-TextValue "D34" "5.525"
Set-TextValue "E34" "  -2.75%  "
Set-TextValue "D35" "11.92"
Set-TextValue "E35" "  -5.36%  "
Set-TextValue "D36" "0.02298"
Set-TextValue "E36" "  -1.50%  "
Set-TextValue "D37" "0.2106"
Set-TextValue "E37" "  -0.36%  "
Set-TextValue "D38" "0.06032"
Set-TextValue "E38" "  -0.76%  "
Set-TextValue "D39" "0.6378"
Set-TextValue "D40" "4.962"
Set-TextValue "E40" "  -1.84%  "
Set-TextValue "D41" "1.203"
Set-TextValue "E41" "  +1.40%  "
Set-TextValue "E42" "  -0.60%  "
Set-TextValue "D43" "7.762"
Set-TextValue "E43" "  -0.82%  "
Set-TextValue "D44" "13.34"
Set-TextValue "E44" "  -1.26%  "
Set-TextValue "D45" "0.5899"
Set-TextValue "E45" "  -1.00%  "
Set-TextValue "D46" "3.711"
Set-TextValue "E46" "  +0.40%  "
Set-TextValue "D47" "122.77"
Set-TextValue "E47" "  -0.47%  "
Set-TextValue "D48" "1.952"
Set-TextValue "E48" "  -3.10%  "
Set-TextValue "D49" "1.150"
Set-TextValue "E49" "  -1.92%  "
Set-TextValue "D50" "0.06855"
Set-TextValue "E50" "  -0.42%  "
Set-TextValue "D51" "72.22"
Set-TextValue "E51" "  -2.56%  "
